$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Make room for the new rows by inserting / deleting blank rows so
# that the existing ("old") rows land on their new target row numbers,
# then fill in the new/changed cell values.
#
# Old layout (rows 259-264):
#   259: ถูกหนึ่ง | ประเทศไทย เพราะประเทศไทยมีตรัง
#   260: ถูกหนึ่ง | ไทย เพราะประเทศไทยมีตรัง
#   261: ถูกหนึ่ง | ไทย
#   262: ถูกหนึ่ง | ตรัง
#   264: ถูกสอง   | ไปฉันเพล
#
# New layout (rows 259-270):
#   259: ตกลง     | ok      (new)
#   260: ตกลง     | k       (new)
#   261: ตกลง     | okay    (new)
#   263: ถูกหนึ่ง | ประเทศไทย เพราะประเทศไทยมีตรัง
#   264: ถูกหนึ่ง | ไทย เพราะประเทศไทยมีตรัง
#   265: ถูกหนึ่ง | ไทย เพราะมีตรัง        (new)
#   266: ถูกหนึ่ง | ประเทศไทย เพราะมีตรัง  (new)
#   267: ถูกหนึ่ง | ไทย
#   268: ถูกหนึ่ง | ตรัง
#   270: ถูกสอง   | ไปฉันเพล
# ------------------------------------------------------------------

# 1) Insert 6 blank rows before row 259. This pushes the old rows
#    259,260,261,262,264 down to 265,266,267,268,270 respectively.
$ws.Range("A259:A264").EntireRow.Insert()

# 2) Rows 263 and 264 are now blank (the gap created by the insert).
#    Delete them so the first two "ถูกหนึ่ง" rows slide back up to
#    263 and 264.
$ws.Range("A263:A264").EntireRow.Delete()

# 3) Insert 2 fresh blank rows before (what is now) row 265 to make
#    room for the two brand-new "ถูกหนึ่ง" rows.
$ws.Range("A265:A266").EntireRow.Insert()

# ------------------------------------------------------------------
# Fill the new "ถูกหนึ่ง" rows first so the new shared strings are
# appended in the same order they appear in the target workbook.
# ------------------------------------------------------------------
$ws.Cells.Item(265, 1).Value = "ถูกหนึ่ง"
$ws.Cells.Item(265, 2).Value = "ไทย เพราะมีตรัง"
$ws.Cells.Item(266, 1).Value = "ถูกหนึ่ง"
$ws.Cells.Item(266, 2).Value = "ประเทศไทย เพราะมีตรัง"

# New "ตกลง" training rows.
$ws.Cells.Item(259, 1).Value = "ตกลง"
$ws.Cells.Item(259, 2).Value = "ok"
$ws.Cells.Item(260, 1).Value = "ตกลง"
$ws.Cells.Item(260, 2).Value = "k"
$ws.Cells.Item(261, 1).Value = "ตกลง"
$ws.Cells.Item(261, 2).Value = "okay"

# ------------------------------------------------------------------
# Restore the view state (scroll position / selection) to match the
# edited workbook as closely as this runtime allows.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 253
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E266").Select()
